$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.693.29"
$ws.Range("E2").Value = "  +5.44%  "

$ws.Range("D3").Value = "2.220.00"
$ws.Range("E3").Value = "  +3.10%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.34"
$ws.Range("E5").Value = "  +1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.95"
$ws.Range("E7").Value = "  -2.56%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.400"
$ws.Range("E9").Value = "  +2.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.98"
$ws.Range("E10").Value = "  +1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +5.92%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "2.550.41"
$ws.Range("E13").Value = "  +3.01%  "

$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.68"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.796"
$ws.Range("E16").Value = "  -0.81%  "

$ws.Range("E17").Value = "  +1.73%  "

$ws.Range("D18").Value = "2.221.33"
$ws.Range("E18").Value = "  +3.72%  "

$ws.Range("D19").Value = "41.561.26"
$ws.Range("E19").Value = "  +5.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.67"
$ws.Range("E20").Value = "  +1.62%  "

$ws.Range("D21").Value = "0.0₃0894"
$ws.Range("E21").Value = "  +5.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.15"
$ws.Range("E23").Value = "  +10.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("E26").Value = "  -1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("E27").Value = "  +1.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.65"
$ws.Range("E28").Value = "  -1.67%  "

$ws.Range("E29").Value = "  +1.62%  "

$ws.Range("E30").Value = "  +1.81%  "

$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.64"
$ws.Range("E32").Value = "  -1.70%  "

$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.92"
$ws.Range("E34").Value = "  +4.94%  "

$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("E37").Value = "  -5.59%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.68"
$ws.Range("E38").Value = "  -2.66%  "

$ws.Range("E39").Value = "  -1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.000244"
$ws.Range("E40").Value = "  +30.83%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.86"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("E43").Value = "  +4.83%  "

$ws.Range("E44").Value = "  +9.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0979"
$ws.Range("E45").Value = "  +6.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.50"
$ws.Range("E46").Value = "  -3.93%  "

$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("D48").Value = "1.464.68"
$ws.Range("E48").Value = "  -3.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.46"
$ws.Range("E49").Value = "  -6.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.81"
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.07"
$ws.Range("E51").Value = "  -1.27%  "
